$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the old "Late" column (N), shifting
# the old N/O/P ("Late" / "Outstanding" heading / "Outstanding") to O/P/Q.
$ws.Columns.Item(14).Insert()

# Give the freshly inserted column N the same display width as column M
# (11 characters), without the bestFit flag the neighbouring columns use.
$ws.Columns.Item(14).ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab and select R6, matching
# the saved view state of the workbook after the edit.
$ws.Activate()
$ws.Range("R6").Select()
